$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 285 (pushes existing rows 285-348 down to 286-349,
# carrying their formatting/values with them, and extends the used range
# to A1:R349).
$ws.Rows("285:285").Insert()

# Populate the newly inserted row 285 with the new data record.
$ws.Cells.Item(285, 1).Value  = 9
$ws.Cells.Item(285, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(285, 3).Value  = "Metropolitana"
$ws.Cells.Item(285, 4).Value  = 44641
$ws.Cells.Item(285, 5).Value  = 13
$ws.Cells.Item(285, 6).Value  = 100112032
$ws.Cells.Item(285, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(285, 8).Value  = "Sin especificar"
$ws.Cells.Item(285, 9).Value  = "Primera"
$ws.Cells.Item(285, 10).Value = 160
$ws.Cells.Item(285, 11).Value = 8000
$ws.Cells.Item(285, 12).Value = 9000
$ws.Cells.Item(285, 13).Value = 8500
$ws.Cells.Item(285, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(285, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(285, 16).Value = 142
$ws.Cells.Item(285, 17).Value = 60
$ws.Cells.Item(285, 18).Value = "Hortaliza"
